$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Chocolate" level labels (column B) ---
# with/without -> no/yes, matching the new row order/content below
$ws.Range("B2").Value = "no"
$ws.Range("B3").Value = "yes"
$ws.Range("B4").Value = "no"
$ws.Range("B5").Value = "yes"

# --- Row 2 (female / no) ---
$ws.Range("C2").Value = 6.4049622804385953
$ws.Range("D2").Value = 6.3543604244993226
$ws.Range("E2").Value = 6.4559670955513253
$ws.Range("F2").Value = 148
$ws.Range("G2").Value = 6.4121717987364857
$ws.Range("H2").Value = 0.30439354425456205
$ws.Range("I2").Value = 0.025020995208869232
$ws.Range("J2").Value = 6.4195765219999998
$ws.Range("K2").Value = 6.2231268034999996
$ws.Range("L2").Value = 6.6217078469999997
$ws.Range("M2").Value = 6.3627244745830076
$ws.Range("N2").Value = 6.4616191228899638

# --- Row 3 (female / yes) ---
$ws.Range("C3").Value = 6.4033185758665621
$ws.Range("D3").Value = 6.3527246666789541
$ws.Range("E3").Value = 6.4543154213974221
$ws.Range("F3").Value = 148
$ws.Range("G3").Value = 6.4008627016621622
$ws.Range("H3").Value = 0.30421548608225651
$ws.Range("I3").Value = 0.025006358917264996
$ws.Range("J3").Value = 6.4042741164999999
$ws.Range("K3").Value = 6.2033664315000001
$ws.Range("L3").Value = 6.6057910829999997
$ws.Range("M3").Value = 6.351444302235671
$ws.Range("N3").Value = 6.4502811010886534

# --- Row 4 (male / no) ---
$ws.Range("C4").Value = 7.7805954198148362
$ws.Range("D4").Value = 7.7055035380197001
$ws.Range("E4").Value = 7.8564190890504317
$ws.Range("F4").Value = 99
$ws.Range("G4").Value = 7.7908388779191915
$ws.Range("H4").Value = 0.40153737786211341
$ws.Range("I4").Value = 0.040356024899145117
$ws.Range("J4").Value = 7.7940105600000003
$ws.Range("K4").Value = 7.4779692312500003
$ws.Range("L4").Value = 8.0341843214999997
$ws.Range("M4").Value = 7.7107536599135038
$ws.Range("N4").Value = 7.8709240959248792

# --- Row 5 (male / yes) ---
$ws.Range("C5").Value = 7.7903587737833098
$ws.Range("D5").Value = 7.7151726641582243
$ws.Range("E5").Value = 7.8662775891204557
$ws.Range("F5").Value = 99
$ws.Range("G5").Value = 7.800643225494948
$ws.Range("H5").Value = 0.40273988165626112
$ws.Range("I5").Value = 0.040476881077756231
$ws.Range("J5").Value = 7.7882664259999999
$ws.Range("K5").Value = 7.4907352527500004
$ws.Range("L5").Value = 8.0477386592499993
$ws.Range("M5").Value = 7.7203181723361309
$ws.Range("N5").Value = 7.8809682786537651
